$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to Text format before assignment so Excel does not
    # auto-convert numeric-looking strings into real numbers (which would
    # drop formatting such as trailing zeros or change precision).
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Restore the original (default) cell style so no stray formatting is left behind.
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "45.930.53"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.581.67"
$ws.Range("E3").Value = "  +8.56%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue "D5" "304.66"
$ws.Range("E5").Value = "  +1.21%  "
Set-TextValue "D6" "99.24"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +4.74%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "0.565"
$ws.Range("E9").Value = "  +11.08%  "
Set-TextValue "D10" "38.09"
$ws.Range("E10").Value = "  +10.25%  "
Set-TextValue "D11" "0.0830"
$ws.Range("E11").Value = "  +5.34%  "
Set-TextValue "D12" "8.00"
$ws.Range("E12").Value = "  +11.96%  "
$ws.Range("D13").Value = "2.976.50"
$ws.Range("E13").Value = "  +8.59%  "
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "2.581.48"
$ws.Range("E15").Value = "  +8.02%  "
$ws.Range("E16").Value = "  +7.49%  "
Set-TextValue "D17" "14.68"
$ws.Range("E17").Value = "  +6.70%  "
$ws.Range("D18").Value = "46.015.00"
$ws.Range("E18").Value = "  +0.11%  "
Set-TextValue "D19" "12.94"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("E20").Value = "  +4.79%  "
Set-TextValue "D21" "6.58"
$ws.Range("E21").Value = "  +8.70%  "
Set-TextValue "D22" "70.39"
$ws.Range("E22").Value = "  +5.32%  "
Set-TextValue "D23" "251.87"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("E24").Value = "  +5.74%  "
$ws.Range("E25").Value = "  +13.03%  "
Set-TextValue "D26" "27.48"
$ws.Range("E26").Value = "  +31.03%  "
Set-TextValue "D27" "1.00"
$ws.Range("E27").Value = "  +0.09%  "
Set-TextValue "D28" "10.31"
$ws.Range("E28").Value = "  +5.70%  "
Set-TextValue "D29" "38.98"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  +2.12%  "
Set-TextValue "D31" "6.07"
$ws.Range("E31").Value = "  +9.47%  "
Set-TextValue "D32" "3.65"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  +4.51%  "
Set-TextValue "D34" "2.25"
$ws.Range("E34").Value = "  +17.50%  "
Set-TextValue "D35" "152.13"
$ws.Range("E35").Value = "  +3.33%  "
Set-TextValue "D36" "0.0823"
$ws.Range("E36").Value = "  +6.39%  "
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  +4.57%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D39" "15.70"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D40" "4.12"
$ws.Range("E40").Value = "  +6.09%  "
$ws.Range("E41").Value = "  +9.55%  "
$ws.Range("E42").Value = "  +6.20%  "
$ws.Range("D43").Value = "2.046.51"
$ws.Range("E43").Value = "  +5.91%  "
Set-TextValue "D44" "19.44"
$ws.Range("E44").Value = "  +35.99%  "
$ws.Range("E45").Value = "  -0.09%  "
Set-TextValue "D46" "90.58"
$ws.Range("E46").Value = "  -1.79%  "
Set-TextValue "D47" "9.17"
$ws.Range("E47").Value = "  +7.53%  "
$ws.Range("E48").Value = "  -1.97%  "
Set-TextValue "D49" "108.15"
$ws.Range("E49").Value = "  +9.15%  "
$ws.Range("D50").Value = "2.835.09"
$ws.Range("E51").Value = "  +6.17%  "
